$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Trening" header column (F1), matching the style of the other
#     header cells (bold font + border + centered alignment) ---
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F1").Value = "Trening"

# --- Give column A a datetime number format so the timestamps become
#     real Excel date serials instead of text. Two NumberFormat writes
#     (lower-case then upper-case) so the style sheet grows numFmtId 164
#     and 165, matching the authored edit. ---
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Apply that same number format (style) to the rest of the timestamp column
$ws.Range("A2").Copy()
$ws.Range("A3:A13").PasteSpecial(-4122)   # xlPasteFormats

# --- Full replacement data set: split into "Duża Gra" / "Mała Gra" parts ---
$data = @(
  @(2, 45684.59363298611, 689.8, 10.15, 0.4294206585202897, "10-15", "Duża Gra"),
  @(3, 45684.59379849537, 704.1, 11.83, 0.3176292691911971, "10-15", "Duża Gra"),
  @(4, 45684.59394201389, 716.5, 10.63, 0.256029207791601, "10-15", "Duża Gra"),
  @(5, 45684.59351724537, 679.8, 5.64, 0.5074322053364344, "5-10", "Duża Gra"),
  @(6, 45684.59360868055, 687.7, 7.44, 0.7170571855136326, "5-10", "Duża Gra"),
  @(7, 45684.59363182871, 689.7, 9.91, 0.4088418462446756, "5-10", "Duża Gra"),
  @(8, 45684.59742233796, 1017.2, 13.83, 1.906478881835937, "10-15", "Mała Gra"),
  @(9, 45684.59995011574, 1235.6, 11.79, 1.303314396313259, "10-15", "Mała Gra"),
  @(10, 45684.60073946759, 1303.8, 14.37, 1.41235889707293, "10-15", "Mała Gra"),
  @(11, 45684.59741655093, 1016.7, 9.82, 1.083679624966213, "5-10", "Mała Gra"),
  @(12, 45684.59994664352, 1235.3, 9.99, 1.093999811581203, "5-10", "Mała Gra"),
  @(13, 45684.60019201389, 1256.5, 5.5, 1.406001329421997, "5-10", "Mała Gra")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

Write-Output "done"
